$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").PrefixCharacter = "'"
$ws.Range("D2").Value = '26.251.98'
$ws.Range("E2").PrefixCharacter = "'"
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").PrefixCharacter = "'"
$ws.Range("D3").Value = '1.589.32'
$ws.Range("E3").PrefixCharacter = "'"
$ws.Range("E3").Value = '  +0.53%  '
$ws.Range("E4").PrefixCharacter = "'"
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").PrefixCharacter = "'"
$ws.Range("D5").Value = '212.60'
$ws.Range("E5").PrefixCharacter = "'"
$ws.Range("E5").Value = '  +1.49%  '
$ws.Range("D6").PrefixCharacter = "'"
$ws.Range("D6").Value = '0.501'
$ws.Range("E6").PrefixCharacter = "'"
$ws.Range("E6").Value = '  +1.00%  '
$ws.Range("E7").PrefixCharacter = "'"
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").PrefixCharacter = "'"
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("E9").PrefixCharacter = "'"
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("D10").PrefixCharacter = "'"
$ws.Range("D10").Value = '19.28'
$ws.Range("E10").PrefixCharacter = "'"
$ws.Range("E10").Value = '  -1.38%  '
$ws.Range("D11").PrefixCharacter = "'"
$ws.Range("D11").Value = '0.0850'
$ws.Range("E11").PrefixCharacter = "'"
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").PrefixCharacter = "'"
$ws.Range("D12").Value = '1.812.56'
$ws.Range("E12").PrefixCharacter = "'"
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").PrefixCharacter = "'"
$ws.Range("D13").Value = '1.591.43'
$ws.Range("E13").PrefixCharacter = "'"
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("D14").PrefixCharacter = "'"
$ws.Range("D14").Value = '4.03'
$ws.Range("E14").PrefixCharacter = "'"
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("E15").PrefixCharacter = "'"
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").PrefixCharacter = "'"
$ws.Range("D16").Value = '64.37'
$ws.Range("E16").PrefixCharacter = "'"
$ws.Range("E16").Value = '  -0.15%  '
$ws.Range("D17").PrefixCharacter = "'"
$ws.Range("D17").Value = '26.274.77'
$ws.Range("E17").PrefixCharacter = "'"
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("E18").PrefixCharacter = "'"
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").PrefixCharacter = "'"
$ws.Range("D19").Value = '7.44'
$ws.Range("E19").PrefixCharacter = "'"
$ws.Range("E19").Value = '  +2.38%  '
$ws.Range("D20").PrefixCharacter = "'"
$ws.Range("D20").Value = '213.12'
$ws.Range("E20").PrefixCharacter = "'"
$ws.Range("E20").Value = '  +2.84%  '
$ws.Range("E21").PrefixCharacter = "'"
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").PrefixCharacter = "'"
$ws.Range("D22").Value = '4.28'
$ws.Range("E22").PrefixCharacter = "'"
$ws.Range("E22").Value = '  +0.50%  '
$ws.Range("E23").PrefixCharacter = "'"
$ws.Range("E23").Value = '  +1.30%  '
$ws.Range("D24").PrefixCharacter = "'"
$ws.Range("D24").Value = '2.15'
$ws.Range("E24").PrefixCharacter = "'"
$ws.Range("E24").Value = '  -2.24%  '
$ws.Range("D25").PrefixCharacter = "'"
$ws.Range("D25").Value = '144.65'
$ws.Range("E25").PrefixCharacter = "'"
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E26").PrefixCharacter = "'"
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("E27").PrefixCharacter = "'"
$ws.Range("E27").Value = '  +0.64%  '
$ws.Range("E28").PrefixCharacter = "'"
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("D29").PrefixCharacter = "'"
$ws.Range("D29").Value = '15.16'
$ws.Range("D30").PrefixCharacter = "'"
$ws.Range("D30").Value = '0.0498'
$ws.Range("E30").PrefixCharacter = "'"
$ws.Range("E30").Value = '  -1.23%  '
$ws.Range("E31").PrefixCharacter = "'"
$ws.Range("E31").Value = '  +1.06%  '
$ws.Range("D32").PrefixCharacter = "'"
$ws.Range("D32").Value = '3.20'
$ws.Range("E32").PrefixCharacter = "'"
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("D33").PrefixCharacter = "'"
$ws.Range("D33").Value = '2.94'
$ws.Range("E33").PrefixCharacter = "'"
$ws.Range("E33").Value = '  -0.36%  '
$ws.Range("D34").PrefixCharacter = "'"
$ws.Range("D34").Value = '1.338.37'
$ws.Range("E34").PrefixCharacter = "'"
$ws.Range("E34").Value = '  +4.96%  '
$ws.Range("D35").PrefixCharacter = "'"
$ws.Range("D35").Value = '2.45'
$ws.Range("E35").PrefixCharacter = "'"
$ws.Range("E35").Value = '  -0.78%  '
$ws.Range("E36").PrefixCharacter = "'"
$ws.Range("E36").Value = '  -0.83%  '
$ws.Range("D37").PrefixCharacter = "'"
$ws.Range("D37").Value = '0.592'
$ws.Range("E37").PrefixCharacter = "'"
$ws.Range("E37").Value = '  -2.88%  '
$ws.Range("E38").PrefixCharacter = "'"
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("D39").PrefixCharacter = "'"
$ws.Range("D39").Value = '0.814'
$ws.Range("E39").PrefixCharacter = "'"
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("E40").PrefixCharacter = "'"
$ws.Range("E40").Value = '  -2.39%  '
$ws.Range("D41").PrefixCharacter = "'"
$ws.Range("D41").Value = '5.71'
$ws.Range("E41").PrefixCharacter = "'"
$ws.Range("E41").Value = '  +3.38%  '
$ws.Range("E43").PrefixCharacter = "'"
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("D44").PrefixCharacter = "'"
$ws.Range("D44").Value = '0.762'
$ws.Range("E44").PrefixCharacter = "'"
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("B45").PrefixCharacter = "'"
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").PrefixCharacter = "'"
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").PrefixCharacter = "'"
$ws.Range("D45").Value = '1.724.45'
$ws.Range("E45").PrefixCharacter = "'"
$ws.Range("E45").Value = '  +0.42%  '
$ws.Range("B46").PrefixCharacter = "'"
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").PrefixCharacter = "'"
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").PrefixCharacter = "'"
$ws.Range("D46").Value = '61.75'
$ws.Range("E46").PrefixCharacter = "'"
$ws.Range("E46").Value = '  -0.99%  '
$ws.Range("D47").PrefixCharacter = "'"
$ws.Range("D47").Value = '86.51'
$ws.Range("E47").PrefixCharacter = "'"
$ws.Range("E47").Value = '  -2.78%  '
$ws.Range("D48").PrefixCharacter = "'"
$ws.Range("D48").Value = '1.49'
$ws.Range("E48").PrefixCharacter = "'"
$ws.Range("E48").Value = '  -3.74%  '
$ws.Range("E49").PrefixCharacter = "'"
$ws.Range("E49").Value = '  -0.57%  '
$ws.Range("D50").PrefixCharacter = "'"
$ws.Range("D50").Value = '0.0977'
$ws.Range("E50").PrefixCharacter = "'"
$ws.Range("E50").Value = '  -2.69%  '
$ws.Range("D51").PrefixCharacter = "'"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").PrefixCharacter = "'"
$ws.Range("E51").Value = '  -0.30%  '
